# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit:
#  - Updated timestamp banner in A1
#  - Updated case statistics for several countries
#  - Re-ordered three pairs/trios of countries (Costa de Marfil / Vietnam / Senegal,
#    Trinidad yTobago / Ruanda, Zambia / Puerto Rico) together with their data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "last refreshed" banner ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 22:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 305934
$ws.Range("C4").Value = 28773
$ws.Range("E4").Value = 282942
$ws.Range("F4").Value = 8073
$ws.Range("G4").Value = 902
$ws.Range("H4").Value = 8306

# --- Row 16: Canada ---
$ws.Range("B16").Value = 13901
$ws.Range("C16").Value = 1526
$ws.Range("D16").Value = 2595
$ws.Range("E16").Value = 11075
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = 231

# --- Row 72: Bosnia y Herzegovina ---
$ws.Range("B72").Value = 624
$ws.Range("C72").Value = 45
$ws.Range("D72").Value = 30
$ws.Range("E72").Value = 573

# --- Rows 99-101: Costa de Marfil moves ahead of Vietnam and Senegal,
#     and gets refreshed figures; Vietnam and Senegal keep their previous
#     figures but shift down one row ---
$ws.Range("A99").Value = "Costa de Marfil"
$ws.Range("B99").Value = 245
$ws.Range("C99").Value = 27
$ws.Range("D99").Value = 25
$ws.Range("E99").Value = 219
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 1

$ws.Range("A100").Value = "Vietnam"
$ws.Range("B100").Value = 240
$ws.Range("C100").Value = 1
$ws.Range("D100").Value = 90
$ws.Range("E100").Value = 150
$ws.Range("F100").Value = 3
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0

$ws.Range("A101").Value = "Senegal"
$ws.Range("B101").Value = 219
$ws.Range("C101").Value = 12
$ws.Range("D101").Value = 72
$ws.Range("E101").Value = 145
$ws.Range("F101").Value = 1
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 2

# --- Rows 123-124: Trinidad yTobago moves ahead of Ruanda and gets
#     refreshed figures; Ruanda keeps its previous figures one row down ---
$ws.Range("A123").Value = "Trinidad yTobago"
$ws.Range("B123").Value = 103
$ws.Range("C123").Value = 5
$ws.Range("D123").Value = 1
$ws.Range("E123").Value = 96
$ws.Range("H123").Value = 6

$ws.Range("A124").Value = "Ruanda"
$ws.Range("B124").Value = 102
$ws.Range("C124").Value = 13
$ws.Range("D124").Value = 0
$ws.Range("E124").Value = 102
$ws.Range("H124").Value = 0

# --- Row 127: Liechtenstein ---
$ws.Range("E127").Value = 76
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 1

# --- Row 135: El Salvador ---
$ws.Range("D135").Value = 2
$ws.Range("E135").Value = 51

# --- Rows 144-145: Zambia moves ahead of Puerto Rico (their totals were
#     tied, only the daily-recovered/death figures differ) ---
$ws.Range("A144").Value = "Zambia"
$ws.Range("D144").Value = 2
$ws.Range("H144").Value = 1

$ws.Range("A145").Value = "Puerto Rico"
$ws.Range("D145").Value = 1
$ws.Range("H145").Value = 2
